# Updates for "Saudi Arabia Division 1" worksheet.
# The source data rows for several fixtures were re-ordered/corrected; this manifests as
# pairs of adjacent (or near-adjacent) data rows having their entire record (columns B:AD -
# id, team names, scores, odds, results, etc.) swapped between the two rows, while the
# leading row-index column (A) stays put.
#
# Commit message: Atualização de bases das ligas, do dia: 15-06-2024 às 21:10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of worksheet row numbers whose B:AD contents must be exchanged.
$rowPairs = @(
    @(2, 3),
    @(35, 36),
    @(38, 39),
    @(59, 60),
    @(115, 116),
    @(125, 126),
    @(149, 150),
    @(172, 173),
    @(181, 182),
    @(218, 219),
    @(238, 239),
    @(291, 293)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
